$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.036089337857949
$ws.Range("D2").Value2 = 1.039309648881025
$ws.Range("E2").Value2 = 1.053419385321489
$ws.Range("F2").Value2 = 1.05928091977823
$ws.Range("I2").Value2 = 1.037403139914726
$ws.Range("J2").Value2 = 1.041199557182038
$ws.Range("K2").Value2 = 1.04209520192853
$ws.Range("L2").Value2 = 1.056165405659929
$ws.Range("M2").Value2 = 1.062010852647899
$ws.Range("N2").Value2 = 1.017655052853154

# Row 3
$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.037068240035294
$ws.Range("D3").Value2 = 1.040051552434897
$ws.Range("E3").Value2 = 1.054709362022901
$ws.Range("F3").Value2 = 1.060648669849264
$ws.Range("I3").Value2 = 1.037642384563777
$ws.Range("J3").Value2 = 1.041822130433109
$ws.Range("K3").Value2 = 1.042647477839179
$ws.Range("L3").Value2 = 1.057267202047382
$ws.Range("M3").Value2 = 1.063191400853644
$ws.Range("N3").Value2 = 1.017864042599253

# Row 4
$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.037701482891451
$ws.Range("D4").Value2 = 1.040531421692203
$ws.Range("E4").Value2 = 1.055544677315297
$ws.Range("F4").Value2 = 1.061534354414416
$ws.Range("I4").Value2 = 1.037795929093145
$ws.Range("J4").Value2 = 1.042224204596821
$ws.Range("K4").Value2 = 1.04300398064436
$ws.Range("L4").Value2 = 1.057980199996508
$ws.Range("M4").Value2 = 1.063955415319629
$ws.Range("N4").Value2 = 1.017998949110545

# Row 5
$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.037967657066368
$ws.Range("D5").Value2 = 1.040733112251284
$ws.Range("E5").Value2 = 1.055895992814217
$ws.Range("F5").Value2 = 1.061906855947168
$ws.Range("I5").Value2 = 1.037860176781898
$ws.Range("J5").Value2 = 1.042393051592916
$ws.Range("K5").Value2 = 1.043153649118083
$ws.Range("L5").Value2 = 1.058279960857732
$ws.Range("M5").Value2 = 1.064276637621561
$ws.Range("N5").Value2 = 1.018055586187207

# Row 6
$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.038012346473947
$ws.Range("D6").Value2 = 1.040766974242961
$ws.Range("E6").Value2 = 1.055954989023101
$ws.Range("F6").Value2 = 1.061969410019797
$ws.Range("I6").Value2 = 1.037870946513197
$ws.Range("J6").Value2 = 1.042421390896799
$ws.Range("K6").Value2 = 1.043178767067627
$ws.Range("L6").Value2 = 1.058330293043157
$ws.Range("M6").Value2 = 1.064330574104392
$ws.Range("N6").Value2 = 1.018065091247094

# Row 7
$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.037705039685218
$ws.Range("D7").Value2 = 1.04053411687461
$ws.Range("E7").Value2 = 1.05554937102249
$ws.Range("F7").Value2 = 1.061539331167914
$ws.Range("I7").Value2 = 1.03779678876241
$ws.Range("J7").Value2 = 1.042226461465601
$ws.Range("K7").Value2 = 1.043005981328827
$ws.Range("L7").Value2 = 1.057984205348776
$ws.Range("M7").Value2 = 1.063959707385309
$ws.Range("N7").Value2 = 1.017999706203053

# Row 8
$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.036420198188244
$ws.Range("D8").Value2 = 1.039560418092276
$ws.Range("E8").Value2 = 1.053855212993155
$ws.Range("F8").Value2 = 1.059743022122398
$ws.Range("I8").Value2 = 1.037484255151097
$ws.Range("J8").Value2 = 1.041410118603275
$ws.Range("K8").Value2 = 1.042282023396732
$ws.Range("L8").Value2 = 1.056537751076314
$ws.Range("M8").Value2 = 1.062409800438624
$ws.Range("N8").Value2 = 1.017725748910384

# Row 9
$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.034154817256346
$ws.Range("D9").Value2 = 1.037843188917965
$ws.Range("E9").Value2 = 1.050874507030245
$ws.Range("F9").Value2 = 1.056582651360197
$ws.Range("I9").Value2 = 1.036923865091811
$ws.Range("J9").Value2 = 1.039965705771661
$ws.Range("K9").Value2 = 1.040999765676715
$ws.Range("L9").Value2 = 1.053989313702043
$ws.Range("M9").Value2 = 1.059679516084175
$ws.Range("N9").Value2 = 1.017240525043884

# Row 10
$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.032643659741064
$ws.Range("D10").Value2 = 1.03669742460331
$ws.Range("E10").Value2 = 1.048890351345536
$ws.Range("F10").Value2 = 1.05447894087683
$ws.Range("I10").Value2 = 1.036543778903068
$ws.Range("J10").Value2 = 1.038998790339216
$ws.Range("K10").Value2 = 1.040140533211052
$ws.Range("L10").Value2 = 1.052290517699641
$ws.Range("M10").Value2 = 1.057859790060707
$ws.Range("N10").Value2 = 1.016915382323346

# Row 11
$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.031989093896061
$ws.Range("D11").Value2 = 1.036201077191979
$ws.Range("E11").Value2 = 1.048031869400178
$ws.Range("F11").Value2 = 1.053568741656891
$ws.Range("I11").Value2 = 1.036377658250555
$ws.Range("J11").Value2 = 1.038579162380044
$ws.Range("K11").Value2 = 1.03976743548073
$ws.Range("L11").Value2 = 1.05155493451037
$ws.Range("M11").Value2 = 1.057071913978466
$ws.Range("N11").Value2 = 1.016774198899066

# Row 12
$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.031745924813914
$ws.Range("D12").Value2 = 1.036016678369952
$ws.Range("E12").Value2 = 1.047713089252632
$ws.Range("F12").Value2 = 1.053230758915578
$ws.Range("I12").Value2 = 1.036315722161957
$ws.Range("J12").Value2 = 1.038423151432315
$ws.Range("K12").Value2 = 1.039628693504289
$ws.Range("L12").Value2 = 1.051281705002516
$ws.Range("M12").Value2 = 1.056779271134341
$ws.Range("N12").Value2 = 1.016721697848346

# Row 13
$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.031798086961234
$ws.Range("D13").Value2 = 1.036056234044344
$ws.Range("E13").Value2 = 1.047781464269154
$ws.Range("F13").Value2 = 1.053303252615949
$ws.Range("I13").Value2 = 1.03632901813809
$ws.Range("J13").Value2 = 1.038456622758071
$ws.Range("K13").Value2 = 1.039658461227122
$ws.Range("L13").Value2 = 1.051340313746674
$ws.Range("M13").Value2 = 1.056842043627163
$ws.Range("N13").Value2 = 1.016732962182317

# Row 14
$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.031968994163675
$ws.Range("D14").Value2 = 1.036185835405026
$ws.Range("E14").Value2 = 1.048005516943344
$ws.Range("F14").Value2 = 1.053540801755259
$ws.Range("I14").Value2 = 1.036372543316251
$ws.Range("J14").Value2 = 1.038566269372947
$ws.Range("K14").Value2 = 1.039755970225483
$ws.Range("L14").Value2 = 1.051532349310187
$ws.Range("M14").Value2 = 1.057047723850377
$ws.Range("N14").Value2 = 1.016769860354186

# Row 15
$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.032074291189362
$ws.Range("D15").Value2 = 1.036265682679257
$ws.Range("E15").Value2 = 1.048143576138458
$ws.Range("F15").Value2 = 1.053687177558872
$ws.Range("I15").Value2 = 1.036399329941745
$ws.Range("J15").Value2 = 1.038633807392474
$ws.Range("K15").Value2 = 1.039816027952381
$ws.Range("L15").Value2 = 1.051650668540459
$ws.Range("M15").Value2 = 1.057174451405737
$ws.Range("N15").Value2 = 1.016792586686436

# Row 16
$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.032687096350342
$ws.Range("D16").Value2 = 1.036730360858519
$ws.Range("E16").Value2 = 1.048947339881724
$ws.Range("F16").Value2 = 1.05453936277715
$ws.Range("I16").Value2 = 1.0365547713211
$ws.Range("J16").Value2 = 1.039026619697103
$ws.Range("K16").Value2 = 1.040165272477959
$ws.Range("L16").Value2 = 1.052339335889498
$ws.Range("M16").Value2 = 1.057912080242216
$ws.Range("N16").Value2 = 1.016924743893191

# Row 17
$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.033071432385404
$ws.Range("D17").Value2 = 1.037021781503955
$ws.Range("E17").Value2 = 1.049451697425255
$ws.Range("F17").Value2 = 1.055074107334013
$ws.Range("I17").Value2 = 1.036651863076159
$ws.Range("J17").Value2 = 1.039272766792538
$ws.Range("K17").Value2 = 1.040384064756554
$ws.Range("L17").Value2 = 1.052771318986596
$ws.Range("M17").Value2 = 1.058374794464949
$ws.Range("N17").Value2 = 1.017007537022225

# Row 18
$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.033295587443783
$ws.Range("D18").Value2 = 1.037191740609784
$ws.Range("E18").Value2 = 1.049745945512504
$ws.Range("F18").Value2 = 1.055386084555873
$ws.Range("I18").Value2 = 1.036708346386347
$ws.Range("J18").Value2 = 1.039416248847245
$ws.Range("K18").Value2 = 1.040511581818568
$ws.Range("L18").Value2 = 1.053023288209213
$ws.Range("M18").Value2 = 1.058644695447848
$ws.Range("N18").Value2 = 1.017055790790343

# Row 19
$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.033372014867075
$ws.Range("D19").Value2 = 1.037249688593983
$ws.Range("E19").Value2 = 1.049846287583688
$ws.Range("F19").Value2 = 1.055492472639187
$ws.Range("I19").Value2 = 1.036727580518512
$ws.Range("J19").Value2 = 1.039465157024141
$ws.Range("K19").Value2 = 1.040555044747299
$ws.Range("L19").Value2 = 1.053109203402621
$ws.Range("M19").Value2 = 1.058736726065055
$ws.Range("N19").Value2 = 1.017072237621289

# Row 20
$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.033030199008308
$ws.Range("D20").Value2 = 1.036990517060989
$ws.Range("E20").Value2 = 1.049397577942293
$ws.Range("F20").Value2 = 1.055016727105714
$ws.Range("I20").Value2 = 1.036641461427003
$ws.Range("J20").Value2 = 1.0392463669916
$ws.Range("K20").Value2 = 1.040360600838656
$ws.Range("L20").Value2 = 1.05272497122212
$ws.Range("M20").Value2 = 1.058325148836901
$ws.Range("N20").Value2 = 1.016998658039873

# Row 21
$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.031918667199367
$ws.Range("D21").Value2 = 1.036147671950848
$ws.Range("E21").Value2 = 1.047939536322266
$ws.Range("F21").Value2 = 1.053470846555379
$ws.Range("I21").Value2 = 1.036359732626069
$ws.Range("J21").Value2 = 1.038533985110858
$ws.Range("K21").Value2 = 1.039727260585787
$ws.Range("L21").Value2 = 1.051475799672483
$ws.Range("M21").Value2 = 1.056987155908903
$ws.Range("N21").Value2 = 1.016758996397369

# Row 22
$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.031219605050613
$ws.Range("D22").Value2 = 1.035617549319045
$ws.Range("E22").Value2 = 1.047023375026882
$ws.Range("F22").Value2 = 1.052499500799463
$ws.Range("I22").Value2 = 1.036181259270041
$ws.Range("J22").Value2 = 1.038085258167812
$ws.Range("K22").Value2 = 1.039328146865662
$ws.Range("L22").Value2 = 1.050690388159935
$ws.Range("M22").Value2 = 1.056145960022433
$ws.Range("N22").Value2 = 1.016607968873618

# Row 23
$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.03159021003248
$ws.Range("D23").Value2 = 1.035898595485177
$ws.Range("E23").Value2 = 1.047508996512603
$ws.Range("F23").Value2 = 1.053014372541334
$ws.Range("I23").Value2 = 1.03627599831768
$ws.Range("J23").Value2 = 1.038323214978652
$ws.Range("K23").Value2 = 1.039539810613406
$ws.Range("L23").Value2 = 1.051106751123513
$ws.Range("M23").Value2 = 1.056591889553308
$ws.Range("N23").Value2 = 1.016688063900821

# Row 24
$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.033048830663218
$ws.Range("D24").Value2 = 1.037004644183765
$ws.Range("E24").Value2 = 1.049422032006051
$ws.Range("F24").Value2 = 1.055042654544747
$ws.Range("I24").Value2 = 1.036646161943819
$ws.Range("J24").Value2 = 1.039258296207493
$ws.Range("K24").Value2 = 1.040371203484433
$ws.Range("L24").Value2 = 1.05274591377887
$ws.Range("M24").Value2 = 1.058347581535171
$ws.Range("N24").Value2 = 1.01700267018738

# Row 25
$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.034740630901951
$ws.Range("D25").Value2 = 1.03828730241986
$ws.Range("E25").Value2 = 1.051644556462878
$ws.Range("F25").Value2 = 1.05739910894708
$ws.Range("I25").Value2 = 1.037069884195212
$ws.Range("J25").Value2 = 1.04033982158551
$ws.Range("K25").Value2 = 1.041332035281523
$ws.Range("L25").Value2 = 1.054648109446786
$ws.Range("M25").Value2 = 1.060385270441289
$ws.Range("N25").Value2 = 1.017366259991485
